$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.037057652392349
$ws.Range("D2").Value = 1.045121814495947
$ws.Range("E2").Value = 1.045557894294894
$ws.Range("F2").Value = 1.055633613646125
$ws.Range("I2").Value = 1.03546521782319
$ws.Range("J2").Value = 1.042162594398317
$ws.Range("K2").Value = 1.047890942755361
$ws.Range("L2").Value = 1.048325798211387
$ws.Range("M2").Value = 1.058373534106447
$ws.Range("N2").Value = 1.017989928588109
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.0378954342692
$ws.Range("D3").Value = 1.045857390176067
$ws.Range("E3").Value = 1.046294259165347
$ws.Range("F3").Value = 1.056431709525344
$ws.Range("I3").Value = 1.035553460310747
$ws.Range("J3").Value = 1.042645159540015
$ws.Range("K3").Value = 1.048438092438525
$ws.Range("L3").Value = 1.048873823318666
$ws.Range("M3").Value = 1.058985149210004
$ws.Range("N3").Value = 1.018151343604258
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.038438247521684
$ws.Range("D4").Value = 1.046334319731062
$ws.Range("E4").Value = 1.046771744513101
$ws.Range("F4").Value = 1.056949197302397
$ws.Range("I4").Value = 1.03560930355971
$ws.Range("J4").Value = 1.042957455766143
$ws.Range("K4").Value = 1.048792427882849
$ws.Range("L4").Value = 1.049228770346915
$ws.Range("M4").Value = 1.059381313029544
$ws.Range("N4").Value = 1.018255765384203
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.038666614966002
$ws.Range("D5").Value = 1.046535049466566
$ws.Range("E5").Value = 1.046972718705152
$ws.Range("F5").Value = 1.057167002367914
$ws.Range("I5").Value = 1.035632478805945
$ws.Range("J5").Value = 1.043088754309687
$ws.Range("K5").Value = 1.048941459099096
$ws.Range("L5").Value = 1.049378069475842
$ws.Range("M5").Value = 1.059547956238535
$ws.Range("N5").Value = 1.018299657919454
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.038704968708628
$ws.Range("D6").Value = 1.046568766217396
$ws.Range("E6").Value = 1.04700647713459
$ws.Range("F6").Value = 1.057203587584054
$ws.Range("I6").Value = 1.035636352346332
$ws.Range("J6").Value = 1.043110800412312
$ws.Range("K6").Value = 1.048966486073389
$ws.Range("L6").Value = 1.049403142067319
$ws.Range("M6").Value = 1.059575941942343
$ws.Range("N6").Value = 1.01830702728027
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.038441298317804
$ws.Range("D7").Value = 1.046337000996043
$ws.Range("E7").Value = 1.046774429002087
$ws.Range("F7").Value = 1.056952106632051
$ws.Range("I7").Value = 1.035609614413056
$ws.Range("J7").Value = 1.042959210149275
$ws.Range("K7").Value = 1.048794418977673
$ws.Range("L7").Value = 1.049230764979155
$ws.Range("M7").Value = 1.059383539349452
$ws.Range("N7").Value = 1.018256351904121
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.037340636126784
$ws.Range("D8").Value = 1.045370205181317
$ws.Range("E8").Value = 1.045806542343276
$ws.Range("F8").Value = 1.055903111829685
$ws.Range("I8").Value = 1.035495299346165
$ws.Range("J8").Value = 1.042325669531941
$ws.Range("K8").Value = 1.048075792949335
$ws.Range("L8").Value = 1.04851093491516
$ws.Range("M8").Value = 1.058580146788824
$ws.Range("N8").Value = 1.018044484321737
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.035406663636096
$ws.Range("D9").Value = 1.043674049234835
$ws.Range("E9").Value = 1.044108810031367
$ws.Range("F9").Value = 1.05406290848233
$ws.Range("I9").Value = 1.035284275583572
$ws.Range("J9").Value = 1.04120969121199
$ws.Range("K9").Value = 1.046811794831005
$ws.Range("L9").Value = 1.047245160761752
$ws.Range("M9").Value = 1.057167667124691
$ws.Range("N9").Value = 1.017670981191619
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.034121176647021
$ws.Range("D10").Value = 1.042548413211778
$ws.Range("E10").Value = 1.04298235586435
$ws.Range("F10").Value = 1.052841784280926
$ws.Range("I10").Value = 1.035137188395464
$ws.Range("J10").Value = 1.040466058866629
$ws.Range("K10").Value = 1.045970777081626
$ws.Range("L10").Value = 1.046403191404583
$ws.Range("M10").Value = 1.056228268987954
$ws.Range("N10").Value = 1.017421898266261
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.033565477023241
$ws.Range("D11").Value = 1.042062243171935
$ws.Range("E11").Value = 1.042495886824287
$ws.Range("F11").Value = 1.052314397197043
$ws.Range("I11").Value = 1.035071987436658
$ws.Range("J11").Value = 1.040144159264432
$ws.Range("K11").Value = 1.045607017519396
$ws.Range("L11").Value = 1.046039074922314
$ws.Range("M11").Value = 1.05582205586956
$ws.Range("N11").Value = 1.017314029851521
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.033359206181357
$ws.Range("D12").Value = 1.041881845964308
$ws.Range("E12").Value = 1.042315386862785
$ws.Range("F12").Value = 1.052118709791997
$ws.Range("I12").Value = 1.035047542522632
$ws.Range("J12").Value = 1.040024607583373
$ws.Range("K12").Value = 1.045471963706923
$ws.Range("L12").Value = 1.04590389681876
$ws.Range("M12").Value = 1.055671255000304
$ws.Range("N12").Value = 1.017273961154409
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.033403445579211
$ws.Range("D13").Value = 1.041920533238144
$ws.Range("E13").Value = 1.042354095801625
$ws.Range("F13").Value = 1.052160675974013
$ws.Range("I13").Value = 1.035052796272827
$ws.Range("J13").Value = 1.040050251070491
$ws.Range("K13").Value = 1.045500930346034
$ws.Range("L13").Value = 1.04593288974379
$ws.Range("M13").Value = 1.055703598434314
$ws.Range("N13").Value = 1.017282556088468
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.033548423729547
$ws.Range("D14").Value = 1.042047327635778
$ws.Range("E14").Value = 1.042480962624947
$ws.Range("F14").Value = 1.052298217366879
$ws.Range("I14").Value = 1.035069971426533
$ws.Range("J14").Value = 1.040134276747215
$ws.Range("K14").Value = 1.04559585264119
$ws.Range("L14").Value = 1.046027899601148
$ws.Range("M14").Value = 1.05580958888318
$ws.Range("N14").Value = 1.017310717789451
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.033637768256616
$ws.Range("D15").Value = 1.042125474820284
$ws.Range("E15").Value = 1.042559155534384
$ws.Range("F15").Value = 1.052382988737036
$ws.Range("I15").Value = 1.035080523628808
$ws.Range("J15").Value = 1.040186049924272
$ws.Range("K15").Value = 1.045654345753936
$ws.Range("L15").Value = 1.046086447761742
$ws.Range("M15").Value = 1.055874904382456
$ws.Range("N15").Value = 1.017328068950978
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.034158076238355
$ws.Range("D16").Value = 1.042580704968101
$ws.Range("E16").Value = 1.04301466862912
$ws.Range("F16").Value = 1.052876814225367
$ws.Range("I16").Value = 1.035141483785646
$ws.Range("J16").Value = 1.040487424471702
$ws.Range("K16").Value = 1.045994927309413
$ws.Range("L16").Value = 1.046427366481719
$ws.Range("M16").Value = 1.05625523980561
$ws.Range("N16").Value = 1.017429056889809
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.034484700542047
$ws.Range("D17").Value = 1.042866591636072
$ws.Range("E17").Value = 1.043300747581355
$ws.Range("F17").Value = 1.05318694561257
$ws.Range("I17").Value = 1.035179318354394
$ws.Range("J17").Value = 1.040676496035393
$ws.Range("K17").Value = 1.04620867511998
$ws.Range("L17").Value = 1.046641340548552
$ws.Range("M17").Value = 1.056493963328725
$ws.Range("N17").Value = 1.017492400609707
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.034675304073311
$ws.Range("D18").Value = 1.04303346378519
$ws.Range("E18").Value = 1.04346773721819
$ws.Range("F18").Value = 1.053367971827996
$ws.Range("I18").Value = 1.035201240760849
$ws.Range("J18").Value = 1.040786787650455
$ws.Range("K18").Value = 1.046333389733231
$ws.Range("L18").Value = 1.046766192452613
$ws.Range("M18").Value = 1.056633259908194
$ws.Range("N18").Value = 1.017529346557423
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.034740310009839
$ws.Range("D19").Value = 1.043090383048008
$ws.Range("E19").Value = 1.043524697445433
$ws.Range("F19").Value = 1.053429719383558
$ws.Range("I19").Value = 1.03520869097457
$ws.Range("J19").Value = 1.040824395787126
$ws.Range("K19").Value = 1.046375920782539
$ws.Range("L19").Value = 1.046808771216735
$ws.Range("M19").Value = 1.056680765409775
$ws.Range("N19").Value = 1.017541943922274
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.034449647616205
$ws.Range("D20").Value = 1.042835906348587
$ws.Range("E20").Value = 1.04327004111208
$ws.Range("F20").Value = 1.053153657799064
$ws.Range("I20").Value = 1.035175274145103
$ws.Range("J20").Value = 1.040656209466427
$ws.Range("K20").Value = 1.046185737943772
$ws.Range("L20").Value = 1.046618378546969
$ws.Range("M20").Value = 1.056468345039795
$ws.Range("N20").Value = 1.017485604566259
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.033505727378059
$ws.Range("D21").Value = 1.042009984640166
$ws.Range("E21").Value = 1.042443598072386
$ws.Range("F21").Value = 1.052257709105675
$ws.Range("I21").Value = 1.035064920015506
$ws.Range("J21").Value = 1.040109532822074
$ws.Range("K21").Value = 1.045567898647979
$ws.Range("L21").Value = 1.045999919594332
$ws.Range("M21").Value = 1.055778374983496
$ws.Range("N21").Value = 1.01730242490947
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.032913062442441
$ws.Range("D22").Value = 1.041491783546782
$ws.Range("E22").Value = 1.041925117268284
$ws.Range("F22").Value = 1.051695593004899
$ws.Range("I22").Value = 1.034994226591952
$ws.Range("J22").Value = 1.039765909741116
$ws.Range("K22").Value = 1.045179802023548
$ws.Range("L22").Value = 1.045611481301854
$ws.Range("M22").Value = 1.05534505474818
$ws.Range("N22").Value = 1.017187243789972
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.033227167462979
$ws.Range("D23").Value = 1.041766387881469
$ws.Range("E23").Value = 1.042199865324473
$ws.Range("F23").Value = 1.051993466694135
$ws.Range("I23").Value = 1.035031826405242
$ws.Range("J23").Value = 1.039948061431428
$ws.Range("K23").Value = 1.045385504405431
$ws.Range("L23").Value = 1.045817360264167
$ws.Range("M23").Value = 1.055574718842275
$ws.Range("N23").Value = 1.017248304141972
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.034465486249473
$ws.Range("D24").Value = 1.042849771340842
$ws.Range("E24").Value = 1.043283915659001
$ws.Range("F24").Value = 1.053168698722096
$ws.Range("I24").Value = 1.03517710199996
$ws.Range("J24").Value = 1.040665376064944
$ws.Range("K24").Value = 1.046196102145675
$ws.Range("L24").Value = 1.046628753949955
$ws.Range("M24").Value = 1.056479920677724
$ws.Range("N24").Value = 1.017488675410068
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.035905974987697
$ws.Range("D25").Value = 1.044111649804533
$ws.Range("E25").Value = 1.044546777119997
$ws.Range("F25").Value = 1.054537653496060
$ws.Range("I25").Value = 1.035339962114484
$ws.Range("J25").Value = 1.041498142212015
$ws.Range("K25").Value = 1.047138285013327
$ws.Range("L25").Value = 1.047572069128958
$ws.Range("M25").Value = 1.057532437058712
$ws.Range("N25").Value = 1.017767557172239
